# Sea Fish Part 1 Public Release
# Applies the "Sea Fish Part 1" content update to the Creatures and
# NewClassFeatures sheets: re-tags a batch of existing creatures as
# "Needs Review", re-points the Cauldron Toad family at the (merged)
# "Toads" document, and adds the seven new Barracuda/Lionfish/Pufferfish
# creatures plus the new source-document bookkeeping rows.

$wb = $excel.ActiveWorkbook
$creatures = $wb.Worksheets.Item("Creatures")
$classFeatures = $wb.Worksheets.Item("NewClassFeatures")

# ---------------------------------------------------------------------
# Creatures!E70:E75 - Toads family: "Needs Clean Up" -> "Needs Review"
# ---------------------------------------------------------------------
$creatures.Range("E70:E75").Value = "Needs Review"

# ---------------------------------------------------------------------
# Creatures!D82:D85 - Cauldron Toads re-pointed from "A Wicked Brew" to
# the "Toads" document, and their review status updated as well.
# ---------------------------------------------------------------------
$creatures.Range("D82:D85").Value = "Toads"
$creatures.Range("E82:E85").Value = "Needs Review"

# Re-link the Cauldron Toad hyperlinks (D82:D85) to point at the Toads
# document instead of the old "A Wicked Brew" document. Remove the
# stale hyperlinks first, then add fresh ones mirroring D70/D71:D75.
$toadsUrl = "https://editor.gmbinder.com/documents/edit/-MmBT4ZMXK4HQN_Zbkgi"
for ($r = 82; $r -le 85; $r++) {
    $cell = $creatures.Cells.Item($r, 4)
    if ($cell.Hyperlinks.Count -gt 0) {
        $cell.Hyperlinks.Delete()
    }
}
$creatures.Hyperlinks.Add($creatures.Cells.Item(82, 4), $toadsUrl) | Out-Null
$creatures.Hyperlinks.Add($creatures.Range("D83:D85"), $toadsUrl, "", "", "Toads") | Out-Null

# ---------------------------------------------------------------------
# Creatures!A110:F116 - seven new Sea Fish creatures.
# Copy the formatting of an existing data row down first so the new
# rows inherit the right number formats / hyperlink font styles.
# ---------------------------------------------------------------------
$creatures.Range("A85:F85").Copy() | Out-Null
$creatures.Range("A110:F116").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Column A first (matches the order new shared strings were created in)
$creatures.Cells.Item(110, 1).Value = "Barracuda"
$creatures.Cells.Item(111, 1).Value = "Giant Barracuda"
$creatures.Cells.Item(112, 1).Value = "School of Barracuda"
$creatures.Cells.Item(113, 1).Value = "Lionfish"
$creatures.Cells.Item(114, 1).Value = "Giant Lionfish"
$creatures.Cells.Item(115, 1).Value = "Pufferfish"
$creatures.Cells.Item(116, 1).Value = "Giant Pufferfish"

# Column B - challenge ratings
$creatures.Cells.Item(110, 2).Value = 0.125
$creatures.Cells.Item(111, 2).Value = 1
$creatures.Cells.Item(112, 2).Value = 3
$creatures.Cells.Item(113, 2).Value = 0
$creatures.Cells.Item(114, 2).Value = 3
$creatures.Cells.Item(115, 2).Value = 0.5
$creatures.Cells.Item(116, 2).Value = 4

# Column C - creature type (all Beasts)
$creatures.Range("C110:C116").Value = "Beast"

# Column D - source document names
$creatures.Cells.Item(110, 4).Value = "Sea Fish Part 1"
$creatures.Cells.Item(111, 4).Value = "Sea Fish Part 2"
$creatures.Cells.Item(112, 4).Value = "Sea Fish Part 3"
$creatures.Cells.Item(113, 4).Value = "Sea Fish Part 4"
$creatures.Cells.Item(114, 4).Value = "Sea Fish Part 5"
$creatures.Cells.Item(115, 4).Value = "Sea Fish Part 6"
$creatures.Cells.Item(116, 4).Value = "Sea Fish Part 7"

# Column E / F - development + release status
$creatures.Range("E110:E116").Value = "Complete"
$creatures.Range("F110:F116").Value = "Publicly Released"

# Hyperlinks for the new Sea Fish documents.
# NOTE: adding a hyperlink to a multi-cell range re-writes the *first*
# cell of that range with the supplied display text, so re-apply the
# real per-row document names afterwards to undo that side effect.
$seaFishPart1Url = "https://editor.gmbinder.com/documents/edit/-N7yaJRL3kLMr593Ja-3"
$creatures.Hyperlinks.Add($creatures.Cells.Item(110, 4), $seaFishPart1Url) | Out-Null
$creatures.Hyperlinks.Add($creatures.Range("D111:D116"), $seaFishPart1Url, "", "", "Sea Fish Part 1") | Out-Null
$creatures.Cells.Item(111, 4).Value = "Sea Fish Part 2"

$creatures.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------
# NewClassFeatures!B53 - "Needs Clean Up" -> "Needs Review"
# ---------------------------------------------------------------------
$classFeatures.Cells.Item(53, 2).Value = "Needs Review"

# ---------------------------------------------------------------------
# NewClassFeatures!A77:E77 - bookkeeping row for the new Sea Fish Part 1
# document (Name / Development Status / Copyright Safe? / Release
# Status / Source Doc).
# ---------------------------------------------------------------------
$classFeatures.Range("A76:F76").Copy() | Out-Null
$classFeatures.Range("A77:F77").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$classFeatures.Cells.Item(77, 1).Value = "Sea Fish Part 1"
$classFeatures.Cells.Item(77, 2).Value = "Complete"
$classFeatures.Cells.Item(77, 3).Value = "Yes"
$classFeatures.Cells.Item(77, 4).Value = "Publicly Released"
$classFeatures.Cells.Item(77, 5).Value = "Source Doc"

$classFeatures.Hyperlinks.Add($classFeatures.Cells.Item(77, 1), $seaFishPart1Url) | Out-Null

$classFeatures.Range("A1").Select() | Out-Null

# Make the Creatures tab the active sheet/tab, matching the saved state.
$creatures.Activate()
